$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that hold numeric-looking text in column D must be forced to
# Text format first, so Excel keeps them as literal strings (matching
# the original inline-string formatting) instead of converting to numbers.
$textFormatCells = @(
    "D4",
    "D5",
    "D6",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D14",
    "D17",
    "D19",
    "D20",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D37",
    "D38",
    "D41",
    "D43",
    "D44",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values
$cellValues = [ordered]@{
    "D2" = "47.133.24"
    "E2" = "  +3.92%  "
    "D3" = "2.494.52"
    "E3" = "  +1.49%  "
    "D4" = "0.999"
    "D5" = "322.89"
    "E5" = "  +0.96%  "
    "D6" = "105.23"
    "E6" = "  +0.68%  "
    "E7" = "  +0.76%  "
    "D8" = "1.00"
    "E8" = "  -0.02%  "
    "D9" = "0.543"
    "E9" = "  +2.01%  "
    "D10" = "37.43"
    "E10" = "  +3.90%  "
    "D11" = "0.0814"
    "E11" = "  +1.05%  "
    "E12" = "  +0.32%  "
    "D13" = "18.29"
    "E13" = "  -1.44%  "
    "D14" = "7.24"
    "E14" = "  +2.48%  "
    "D15" = "2.878.84"
    "E15" = "  +1.39%  "
    "D16" = "2.487.89"
    "E16" = "  +1.21%  "
    "D17" = "0.844"
    "E17" = "  +0.15%  "
    "D18" = "47.035.71"
    "E18" = "  +3.94%  "
    "D19" = "12.63"
    "E19" = "  +1.90%  "
    "D20" = "6.56"
    "E20" = "  +2.47%  "
    "D21" = "0.0₃0933"
    "E21" = "  +0.24%  "
    "D22" = "70.96"
    "E22" = "  +2.44%  "
    "D23" = "250.93"
    "E23" = "  +2.65%  "
    "D24" = "2.35"
    "E24" = "  +2.18%  "
    "D25" = "2.54"
    "D26" = "26.19"
    "E26" = "  +2.50%  "
    "E27" = "  -0.13%  "
    "D28" = "10.18"
    "E28" = "  +6.01%  "
    "D29" = "2.25"
    "E29" = "  +2.67%  "
    "D30" = "35.34"
    "E30" = "  +3.82%  "
    "D31" = "0.135"
    "E31" = "  +4.41%  "
    "D32" = "49.55"
    "E32" = "  -0.15%  "
    "D33" = "19.72"
    "E33" = "  -3.58%  "
    "D34" = "5.37"
    "E34" = "  +2.12%  "
    "D35" = "0.0785"
    "E35" = "  +2.27%  "
    "D37" = "1.93"
    "E37" = "  +0.43%  "
    "D38" = "4.61"
    "E38" = "  +1.51%  "
    "E39" = "  +2.73%  "
    "E40" = "  +1.09%  "
    "D41" = "121.80"
    "E41" = "  -2.79%  "
    "E42" = "  +1.22%  "
    "D43" = "21.64"
    "E43" = "  +1.36%  "
    "D44" = "0.0295"
    "E44" = "  +1.03%  "
    "D45" = "1.952.46"
    "E45" = "  +0.05%  "
    "E46" = "  -0.42%  "
    "D47" = "2.10"
    "E47" = "  -0.95%  "
    "B48" = "FraxShare"
    "C48" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
    "D48" = "9.15"
    "E48" = "  -1.43%  "
    "B49" = "Stacks"
    "C49" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D49" = "1.79"
    "E49" = "  -0.06%  "
    "D50" = "5.36"
    "E50" = "  +12.71%  "
    "D51" = "78.80"
    "E51" = "  +3.38%  "
}
foreach ($cellRef in $cellValues.Keys) {
    $ws.Range($cellRef).Value = $cellValues[$cellRef]
}
